$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.858.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5021"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2574"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06424"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07704"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.248"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.632.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.860.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5453"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.874.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.332"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.958"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.987"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.930"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1146"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.715"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.242"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05011"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.269"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.361"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.173.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8953"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.605"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5623"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01561"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -3.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.771.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4516"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05055"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
